# Add change password, editing users
#
# - Remove the "Salo" worksheet entirely (it was an empty, unused sheet).
# - Rework the shared-string text used in row 1 of "Jkhfhg":
#     A1: "Saf"        -> "Good"
#     E1: "Maplehgjhf" -> "Paint"

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the empty "Salo" sheet (sheetId 2). Excel renumbers the remaining
# sheets' relationship ids automatically, so "Fender" keeps sheetId 3 but
# moves from rId3 to rId2.
$salo = $wb.Worksheets.Item("Salo")
$salo.Delete()

# Update the two changed labels on the first sheet.
$ws = $wb.Worksheets.Item("Jkhfhg")
$ws.Range("A1").Value = "Good"
$ws.Range("E1").Value = "Paint"
